# Developed Test Cases - Learning & Development
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing data row (row 2) ---
# buyIn1 for TEST-0608 changes denom-count from "100" to "100-3"
$ws.Range("B2").Value = "100-3;rated-6009;1"

# --- New test case: TEST-18052 (row 3) ---
$ws.Range("A3").Value = "TEST-18052"
$ws.Range("B3").Value = "25;anon"
$ws.Range("I3").Value = "P1;25;P1"
$ws.Range("P3").Value = "4d"
$ws.Range("Q3").Value = "2s"
$ws.Range("R3").Value = "4d"
$ws.Range("S3").Value = "3s"
$ws.Range("W3").Value = "P1;25"

# --- New test case: TEST-28843 (row 4) ---
$ws.Range("A4").Value = "TEST-28843"

# --- New test case: TEST-14009 (row 5) ---
$ws.Range("A5").Value = "TEST-14009"
$ws.Range("B5").Value = "1000;anon"
$ws.Range("I5").Value = "P1;1000;B3"
$ws.Range("P5").Value = "Ah"
$ws.Range("Q5").Value = "4s"
$ws.Range("R5").Value = "2s"
$ws.Range("S5").Value = "Ah"
$ws.Range("V5").Value = "B3"

# --- View state: scroll / active cell to match author's last position ---
$excel.ActiveWindow.ScrollColumn = 13
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("W7").Select() | Out-Null
